# "the guarny e lobby 1 winner"
# Lobby roster shuffle: TheGuarny moves from lobby1 -> lobby3 (row 40 of the
# external LOBBY cache), RMT _SILVERMAN-AC moves from lobby2 -> lobby1, and
# the remaining lobby2 (column F/G) rows shift up by one slot.
#
# These cells hold formulas pointing at an external, unavailable workbook
# ([1]LOBBY!... / [1]nomi!...), so the cached <v> is the only thing that can
# actually change here (there is nothing to recalculate against). We set the
# literal display values directly.

$wb = $excel.ActiveWorkbook

$lobby1 = $wb.Worksheets.Item("lobby1")
$lobby2 = $wb.Worksheets.Item("lobby2")
$lobby3 = $wb.Worksheets.Item("lobby3")

# --- lobby1: row 12 -> RMT _SILVERMAN-AC takes TheGuarny's old seat ---
$lobby1.Range("B12").Value = "RMT _SILVERMAN-AC"
$lobby1.Range("C12").Value = "L.Silvestri"
$lobby1.Range("D12").Value = "RMT"

# --- lobby2: rows 10-13 shift up one (Silverman leaves), row 14 empties ---
$lobby2.Range("B10").Value = "Talsigiano"
$lobby2.Range("C10").Value = "Talsigiano "
$lobby2.Range("D10").Value = "TLMA"

$lobby2.Range("B11").Value = "jack-187-jack"
$lobby2.Range("C11").Value = "Jack187"
$lobby2.Range("D11").Value = "TLMA"

$lobby2.Range("B12").Value = "TLM_Verce90"
$lobby2.Range("C12").Value = "Verce"
$lobby2.Range("D12").Value = "TLM"

$lobby2.Range("B13").Value = "Mancinelli2002"
$lobby2.Range("C13").Value = "Gian"
$lobby2.Range("D13").Value = "TLMA"

$lobby2.Range("B14").Value = 0
$lobby2.Range("C14").Value = 0

# --- lobby3: row 15 -> TheGuarny is promoted here (lobby 1 winner) ---
$lobby3.Range("B15").Value = "TheGuarny"
$lobby3.Range("C15").Value = "TLM_TheGuarny"
$lobby3.Range("D15").Value = "TLM"
